$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 ("io" / "E") — shifts rows 6-7 up.
$ws.Rows("5").Delete()

# Match the resulting selection left behind in the file.
$ws.Range("I18").Select()
